$wb = $excel.ActiveWorkbook

# Map of new "想去人数" (F column) values for each affected row,
# applied identically to both the "展览" and "全部类型" sheets.
$updates = @{
    4  = 13149
    5  = 1345
    6  = 224
    9  = 167
    12 = 4
    16 = 48
    17 = 424
    18 = 5564
    22 = 16
    25 = 156
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
